# Daily attendance processing - reorder "Recorded By" (column G) names.
# For each data row, the comma-separated list of recorders is reversed
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"),
# except for entries that include "admin@admin.com", which are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        if ($val -notlike "*admin@admin.com*") {
            $parts = $val -split ", "
            if ($parts.Length -gt 1) {
                $reversed = $parts[($parts.Length - 1)..0]
                $cell.Value2 = ($reversed -join ", ")
            }
        }
    }
}
